$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new column before column N (14), shifting Late/Outstanding etc. one to the right
$ws.Columns("N").Insert()

# New column approximately matches the width of column M (the closest the engine can snap to)
$ws.Columns("N").ColumnWidth = 9.75

# Make "Repayment schedule" the active sheet/tab, with S8 selected
$ws.Activate()
$ws.Range("S8").Select()
